$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "24-11-2025"
$ws.Range("B69").Value = "The price of gold in India today is ₹12,513 per gram for 24 karat gold, ₹11,470 per gram for 22 karat gold and ₹9,385 per gram for 18 karat gold (also called 999 gold)."
